# Apply targeted updates to column F (dSF) values as described in the diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F5").Value = -6
$ws.Range("F7").Value = 2
$ws.Range("F8").Value = -4
$ws.Range("F14").Value = 1
$ws.Range("F15").Value = 1
$ws.Range("F16").Value = 5
$ws.Range("F18").Value = 1
$ws.Range("F19").Value = -2
